$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.951.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.43%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.924.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.38%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.02%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.03%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4776"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.37%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2894"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.87%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06718"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.55%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.92"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.34%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "102.73"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.88%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07833"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.19%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.938.31"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.03%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.247"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.57%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6917"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.83%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "291.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +9.75%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "31.031.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.45%  "

# Row 18
$ws.Range("B18").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C18").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.190.49"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.79%  "

# Row 19
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9999"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.32%  "

# Row 20
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.61%  "

# Row 21
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000007548"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.26%  "

# Row 22
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.519"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.88%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.005"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.14%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.365"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.67%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.471"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.44%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.78%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.17%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.069"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.19%  "

# Row 29
$ws.Range("E29").Value = "  +0.27%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.1001"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.96%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.576"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.85%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.532"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.22%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.297"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.26%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04809"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.27%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7310"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.64%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.122"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.17%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.729"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.71%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01946"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.52%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.835"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.05%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.623"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.64%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "75.54"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.86%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.015"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.56%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8720"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.55%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4327"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.12%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "105.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.82%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.002"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.05%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.021.76"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.98%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.486"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.07%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.269"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.98%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1198"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.36%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.97"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.45%  "
